# Auto-generated edit script applying the Sargatanas_Profits.xlsx market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 971.6786
$ws.Range("J17").Value = 971.6786
$ws.Range("L17").Value = 2915.0358
$ws.Range("N17").Value = -3251.0358

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H116").Value = 41674500
$ws.Range("I116").Value = 83335660
$ws.Range("K116").Value = 83335660
$ws.Range("M116").Value = -83332218

$ws.Range("H132").Value = 1647.5151
$ws.Range("I132").Value = 1633.8276
$ws.Range("J132").Value = 1746.75
$ws.Range("K132").Value = 4901.4828
$ws.Range("L132").Value = 5240.25
$ws.Range("M132").Value = -2371.4828
$ws.Range("N132").Value = -10300.25

$ws.Range("H137").Value = 1596.5217
$ws.Range("I137").Value = 1385.4736
$ws.Range("K137").Value = 4156.4208
$ws.Range("M137").Value = -1606.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2543.9092
$ws.Range("I2").Value = 799.6429000000001
$ws.Range("K2").Value = 799.6429000000001
$ws.Range("M2").Value = -686.6429000000001

$ws.Range("H32").Value = 3129385.8
$ws.Range("I32").Value = 3208344
$ws.Range("K32").Value = 3208344
$ws.Range("M32").Value = -3208057

$ws.Range("H61").Value = 8809.77
$ws.Range("I61").Value = 2954.3
$ws.Range("K61").Value = 2954.3
$ws.Range("M61").Value = -2742.3

$ws.Range("H74").Value = 26729.773
$ws.Range("I74").Value = 40857.73
$ws.Range("J74").Value = 6322.722
$ws.Range("K74").Value = 40857.73
$ws.Range("L74").Value = 6322.722
$ws.Range("M74").Value = -39983.73
$ws.Range("N74").Value = -8070.722

$ws.Range("H77").Value = 26729.773
$ws.Range("I77").Value = 40857.73
$ws.Range("J77").Value = 6322.722
$ws.Range("K77").Value = 204288.65
$ws.Range("L77").Value = 31613.61
$ws.Range("M77").Value = -199920.65
$ws.Range("N77").Value = -40349.61

$ws.Range("H88").Value = 2999.5
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 2999.5
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H95").Value = 55000
$ws.Range("J95").Value = 55000
$ws.Range("L95").Value = 55000
$ws.Range("N95").Value = -60492

$ws.Range("H116").Value = 2543.9092
$ws.Range("I116").Value = 799.6429000000001
$ws.Range("K116").Value = 799.6429000000001
$ws.Range("M116").Value = 1494.3571

$ws.Range("H122").Value = 5821.0835
$ws.Range("I122").Value = 4539.222
$ws.Range("K122").Value = 13617.666
$ws.Range("M122").Value = -11167.666

$ws.Range("H132").Value = 948865.7
$ws.Range("I132").Value = 1474353.8
$ws.Range("K132").Value = 4423061.4
$ws.Range("M132").Value = -4420531.4

$ws.Range("H136").Value = 8809.77
$ws.Range("I136").Value = 2954.3
$ws.Range("K136").Value = 8862.900000000001
$ws.Range("M136").Value = -6312.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2543.9092
$ws.Range("I3").Value = 799.6429000000001
$ws.Range("K3").Value = 799.6429000000001
$ws.Range("M3").Value = -685.6429000000001

$ws.Range("H82").Value = 1549.6666
$ws.Range("I82").Value = 1549.6666
$ws.Range("K82").Value = 1549.6666
$ws.Range("M82").Value = -1166.6666

$ws.Range("H85").Value = 1549.6666
$ws.Range("I85").Value = 1549.6666
$ws.Range("K85").Value = 1549.6666
$ws.Range("M85").Value = -223.6666

$ws.Range("H134").Value = 5594.853
$ws.Range("I134").Value = 1639.65
$ws.Range("K134").Value = 4918.950000000001
$ws.Range("M134").Value = -2383.950000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 8932662
$ws.Range("I105").Value = 14286270
$ws.Range("K105").Value = 14286270
$ws.Range("M105").Value = -14284523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42730156
$ws.Range("I4").Value = 19820220
$ws.Range("K4").Value = 59460660
$ws.Range("M4").Value = -59460548

$ws.Range("H5").Value = 2094.9678
$ws.Range("I5").Value = 1859.6
$ws.Range("J5").Value = 2522.9092
$ws.Range("K5").Value = 5578.799999999999
$ws.Range("L5").Value = 7568.7276
$ws.Range("M5").Value = -5466.799999999999
$ws.Range("N5").Value = -7792.7276

$ws.Range("H23").Value = 183.41667
$ws.Range("J23").Value = 196.33333
$ws.Range("L23").Value = 588.99999
$ws.Range("N23").Value = -1058.99999

$ws.Range("H107").Value = 1512.4062
$ws.Range("J107").Value = 1964.85
$ws.Range("L107").Value = 5894.549999999999
$ws.Range("N107").Value = -9734.549999999999

$ws.Range("H127").Value = 1296.3334
$ws.Range("J127").Value = 1296.3334
$ws.Range("L127").Value = 3889.0002
$ws.Range("N127").Value = -13809.0002

$ws.Range("H132").Value = 4478.6665
$ws.Range("I132").Value = 1713.2858
$ws.Range("J132").Value = 6898.375
$ws.Range("K132").Value = 15419.5722
$ws.Range("L132").Value = 62085.375
$ws.Range("M132").Value = -12889.5722
$ws.Range("N132").Value = -67145.375

$ws.Range("H135").Value = 2094.9678
$ws.Range("I135").Value = 1859.6
$ws.Range("J135").Value = 2522.9092
$ws.Range("K135").Value = 16736.4
$ws.Range("L135").Value = 22706.1828
$ws.Range("M135").Value = -14201.4
$ws.Range("N135").Value = -27776.1828

$ws.Range("H139").Value = 190100.69
$ws.Range("I139").Value = 202107.4
$ws.Range("K139").Value = 606322.2
$ws.Range("M139").Value = -601182.2

$ws.Range("H140").Value = 167918.3
$ws.Range("I140").Value = 200765.45
$ws.Range("K140").Value = 602296.3500000001
$ws.Range("M140").Value = -597116.3500000001

$ws.Range("H141").Value = 1757.6364
$ws.Range("I141").Value = 1757.6364
$ws.Range("K141").Value = 5272.9092
$ws.Range("M141").Value = -92.90920000000006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3457.4211
$ws.Range("I132").Value = 2088.84
$ws.Range("K132").Value = 6266.52
$ws.Range("M132").Value = -3736.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 59999
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H55").Value = 847.5
$ws.Range("I55").Value = 879.7619
$ws.Range("J55").Value = 795.38464
$ws.Range("K55").Value = 879.7619
$ws.Range("L55").Value = 795.38464
$ws.Range("M55").Value = -706.7619
$ws.Range("N55").Value = -1141.38464

$ws.Range("H82").Value = 2014867.4
$ws.Range("I82").Value = 2349012.2
$ws.Range("J82").Value = 9998
$ws.Range("K82").Value = 2349012.2
$ws.Range("L82").Value = 9998
$ws.Range("M82").Value = -2348651.2
$ws.Range("N82").Value = -10720

$ws.Range("H85").Value = 2014867.4
$ws.Range("I85").Value = 2349012.2
$ws.Range("J85").Value = 9998
$ws.Range("K85").Value = 2349012.2
$ws.Range("L85").Value = 9998
$ws.Range("M85").Value = -2347764.2
$ws.Range("N85").Value = -12494

$ws.Range("H100").Value = 3225.25
$ws.Range("I100").Value = 2666.3333
$ws.Range("J100").Value = 4902
$ws.Range("K100").Value = 2666.3333
$ws.Range("L100").Value = 4902
$ws.Range("M100").Value = -2125.3333
$ws.Range("N100").Value = -5984

$ws.Range("H122").Value = 9278.772000000001
$ws.Range("I122").Value = 13892.25
$ws.Range("J122").Value = 6642.5
$ws.Range("K122").Value = 41676.75
$ws.Range("L122").Value = 19927.5
$ws.Range("M122").Value = -39226.75
$ws.Range("N122").Value = -24827.5

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 6150.311
$ws.Range("I132").Value = 3732.3809
$ws.Range("J132").Value = 8266
$ws.Range("K132").Value = 11197.1427
$ws.Range("L132").Value = 24798
$ws.Range("M132").Value = -8667.1427
$ws.Range("N132").Value = -29858

$ws.Range("H136").Value = 8076.1455
$ws.Range("I136").Value = 5518.9697
$ws.Range("J136").Value = 11911.909
$ws.Range("K136").Value = 16556.9091
$ws.Range("L136").Value = 35735.727
$ws.Range("M136").Value = -14006.9091
$ws.Range("N136").Value = -40835.727

$ws.Range("H137").Value = 91234
$ws.Range("J137").Value = 91234
$ws.Range("L137").Value = 91234
$ws.Range("N137").Value = -101434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 38193.5
$ws.Range("J101").Value = 38193.5
$ws.Range("L101").Value = 38193.5
$ws.Range("N101").Value = -44683.5

$ws.Range("H107").Value = 763
$ws.Range("J107").Value = 821
$ws.Range("L107").Value = 2463
$ws.Range("N107").Value = -6303
